# Ajout d'une nouvelle journée de présences (colonne BW) :
# - nouvel en-tête de date en BW1 (2025-11-05, numéro de série 45966)
# - pour chaque joueur (lignes 2 à 29), report du statut de présence du jour
#   (les formules COUNTA/COUNTIF des colonnes B:J se recalculent automatiquement
#   puisqu'elles couvrent déjà des plages qui incluent la nouvelle colonne)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- En-tête : nouvelle date (même style que BV1) --
$ws.Range("BW1").Value2 = 45966
$ws.Range("BV1").Copy()
$ws.Range("BW1").PasteSpecial(-4122)

# -- Statuts journaliers par joueur (même style que la cellule BV correspondante) --
$dayStatus = @{
    2  = "P"
    3  = "R"
    4  = "P"
    5  = "B"
    6  = "B"
    7  = "P"
    8  = "RH"
    9  = "P"
    10 = "P"
    11 = "P"
    13 = "B"
    14 = "P"
    15 = "P"
    16 = "P"
    17 = "P"
    18 = "B"
    19 = "P"
    20 = "B"
    22 = "P"
    23 = "RH"
    24 = "P"
    25 = "P"
    26 = "P"
    27 = "P"
    28 = "P"
    29 = "B"
}

foreach ($row in $dayStatus.Keys) {
    $src = $ws.Range("BV$row")
    $dst = $ws.Range("BW$row")
    $dst.Value2 = $dayStatus[$row]
    $src.Copy()
    $dst.PasteSpecial(-4122)
}

# -- Ligne 12 : le joueur n'a plus de suivi après la colonne AX, pas de BW12 --
# -- Ligne 21 : cellule vide (même mise en forme que BV21, sans valeur) --
$ws.Range("BV21").Copy()
$ws.Range("BW21").PasteSpecial(-4122)

# -- Sélection finale affichée dans le classeur --
$ws.Range("BY28").Select() | Out-Null
